$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "IRS-Bilateral"

$ws.Range("A1").Value2 = "Value Date"
$ws.Range("B1").Value2 = "Position Account ID"
$ws.Range("C1").Value2 = "Client ID"
$ws.Range("D1").Value2 = "UTI"
$ws.Range("E1").Value2 = "Currency"
$ws.Range("F1").Value2 = "Effective Date"
$ws.Range("G1").Value2 = "Maturity Date"
$ws.Range("H1").Value2 = "Cleared Date"
$ws.Range("I1").Value2 = "Trade type"
$ws.Range("J1").Value2 = "Firm ID"
$ws.Range("K1").Value2 = "Source"
$ws.Range("L1").Value2 = "LEG1_TYPE"
$ws.Range("M1").Value2 = "LEG1_CCY"
$ws.Range("N1").Value2 = "LEG1_PAY_FREQ"
$ws.Range("O1").Value2 = "LEG1_PAY_ADJ_BUS_DAY_CONV"
$ws.Range("P1").Value2 = "LEG1_PAY_ADJ_CAL"
$ws.Range("Q1").Value2 = "LEG1_DAYCOUNT"
$ws.Range("R1").Value2 = "LEG1_INDEX"
$ws.Range("S1").Value2 = "LEG1_INDEX_TENOR"
$ws.Range("T1").Value2 = "LEG1_RESET_FREQ"
$ws.Range("U1").Value2 = "LEG1_START_DATE"
$ws.Range("V1").Value2 = "LEG1_MAT_DATE"
$ws.Range("W1").Value2 = "LEG1_NOTIONAL"
$ws.Range("X1").Value2 = "LEG1_FIXED_RATE"
$ws.Range("Y1").Value2 = "LEG2_TYPE"
$ws.Range("Z1").Value2 = "LEG2_CCY"
$ws.Range("AA1").Value2 = "LEG2_PAY_FREQ"
$ws.Range("AB1").Value2 = "LEG2_PAY_ADJ_BUS_DAY_CONV"
$ws.Range("AC1").Value2 = "LEG2_PAY_ADJ_CAL"
$ws.Range("AD1").Value2 = "LEG2_DAYCOUNT"
$ws.Range("AE1").Value2 = "LEG2_INDEX"
$ws.Range("AF1").Value2 = "LEG2_INDEX_TENOR"
$ws.Range("AG1").Value2 = "LEG2_RESET_FREQ"
$ws.Range("AH1").Value2 = "LEG2_START_DATE"
$ws.Range("AI1").Value2 = "LEG2_MAT_DATE"
$ws.Range("AJ1").Value2 = "LEG2_NOTIONAL"
$ws.Range("AK1").Value2 = "LEG2_FIXED_RATE"
$ws.Range("AL1").Value2 = "LEG1_DIRECTION"
$ws.Range("AM1").Value2 = "LEG2_DIRECTION"
$ws.Range("AN1").Value2 = "Counterpart ID"
$ws.Range("AO1").Value2 = "Agreement ID"
$ws.Range("AP1").Value2 = "Jurisdiction"
$c = $ws.Range("A2"); $c.NumberFormat = "DD/MM/YY"; $c.Value2 = 41631
$ws.Range("B2").Value2 = "acc1"
$ws.Range("C2").Value2 = 11811152
$ws.Range("D2").Value2 = 455820
$ws.Range("E2").Value2 = "SGD"
$c = $ws.Range("F2"); $c.NumberFormat = "DD/MM/YY"; $c.Value2 = 41607
$c = $ws.Range("G2"); $c.NumberFormat = "DD/MM/YY"; $c.Value2 = 44164
$c = $ws.Range("H2"); $c.NumberFormat = "DD/MM/YY"; $c.Value2 = 41605
$ws.Range("I2").Value2 = "Bilateral"
$ws.Range("J2").Value2 = 999
$ws.Range("K2").Value2 = "MARKIT_WIRE"
$ws.Range("L2").Value2 = "FIXED"
$ws.Range("M2").Value2 = "SGD"
$ws.Range("N2").Value2 = "6M"
$ws.Range("O2").Value2 = "ModifiedFollowing"
$ws.Range("P2").Value2 = "SGSI"
$ws.Range("Q2").Value2 = "Act/365F"
$c = $ws.Range("U2"); $c.NumberFormat = "DD/MM/YY"; $c.Value2 = 41607
$c = $ws.Range("V2"); $c.NumberFormat = "DD/MM/YY"; $c.Value2 = 44164
$c = $ws.Range("W2"); $c.NumberFormat = "@"; $c.Value2 = "10,000,000.00"; $c.NumberFormat = "General"
$c = $ws.Range("X2"); $c.NumberFormat = "@"; $c.Value2 = "1.1"; $c.NumberFormat = "General"
$ws.Range("Y2").Value2 = "FLOAT"
$ws.Range("Z2").Value2 = "SGD"
$ws.Range("AA2").Value2 = "6M"
$ws.Range("AB2").Value2 = "ModifiedFollowing"
$ws.Range("AC2").Value2 = "SGSI"
$ws.Range("AD2").Value2 = "Act/365F"
$ws.Range("AE2").Value2 = "SGD-SOR-Reuters"
$ws.Range("AF2").Value2 = "6M"
$ws.Range("AG2").Value2 = "6M"
$c = $ws.Range("AH2"); $c.NumberFormat = "DD/MM/YY"; $c.Value2 = 41607
$c = $ws.Range("AI2"); $c.NumberFormat = "DD/MM/YY"; $c.Value2 = 44164
$c = $ws.Range("AJ2"); $c.NumberFormat = "@"; $c.Value2 = "10,000,000.00"; $c.NumberFormat = "General"
$ws.Range("AL2").Value2 = "R"
$ws.Range("AM2").Value2 = "P"
$ws.Range("AN2").Value2 = 11911171
$ws.Range("AO2").Value2 = 12011171
$ws.Range("AP2").Value2 = "Singapore"

[void]$ws.Range("A3").Select()
$ws.Activate()
